$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.059.06'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = '1.598.78'
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("E4").Value = '  -0.41%  '

$ws.Range("E5").Value = '  -0.38%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '302.32'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("E7").Value = '  +0.14%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3635'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.83%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '50.92'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +4.44%  '

$ws.Range("E10").Value = '  -2.29%  '

$ws.Range("E11").Value = '  -0.34%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.08134'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.41%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '22.31'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.71%  '

$ws.Range("E14").Value = '  -1.13%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '7.334'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.17%  '

$ws.Range("E16").Value = '  -1.82%  '

$ws.Range("D17").Value = '1.597.24'
$ws.Range("E17").Value = '  +0.16%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '92.25'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.01%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06845'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.30%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '18.14'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.83%  '

$ws.Range("E21").Value = '  -1.72%  '

$ws.Range("E22").Value = '  -0.30%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '12.99'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.09%  '

$ws.Range("D24").Value = '23.063.53'
$ws.Range("E24").Value = '  -0.43%  '

$ws.Range("E25").Value = '  +0.56%  '

$ws.Range("E26").Value = '  -4.60%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '21.08'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.36%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '148.87'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.51%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.259'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.60%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '134.60'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.73%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '2.381'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.39%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.722'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -5.16%  '

$ws.Range("D33").Value = '1.776.36'
$ws.Range("E33").Value = '  +0.32%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.9582'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.77%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.07492'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.54%  '

$ws.Range("E36").Value = '  -2.61%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '10.18'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.25%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '6.195'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.24%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.2514'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.36%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.08819'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.58%  '

$ws.Range("E41").Value = '  -2.36%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.7031'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.79%  '

$ws.Range("E43").Value = '  -3.40%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '15.22'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -5.63%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.6579'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.91%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '4.016'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.08%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.269'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.12%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '132.02'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.75%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.07926'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.76%  '

$ws.Range("E50").Value = '  +3.92%  '

$ws.Range("E51").Value = '  +3.81%  '
